$d = $word.ActiveDocument

# --- Merge split runs (no visible text change) by re-finding/replacing the
# --- full sentences so Word recombines the runs into a single run, matching
# --- the target canonical OOXML (which has each paragraph's text in one run). ---
$d.Content.Find.Execute(
    "This week during one of the classes we played a game in which one player used one plastic circle and the other player used 4. The game was that there was a board and the player with 4 circles had to get the 4 pieces to the other side of the board; the same for the player with one piece. I could tell that it seemed that it was easier for the one-piece player to win, however, once the player with 4 pieces learned how to play with the 4 figures, he had a better chance of winning. The way I plan to contribute to my team is by keeping a positive attitude, not being so introverted, but more extroverted so that I can give my suggestions. Likewise, I plan to contribute to the team, helping them to maintain a team atmosphere so we can accomplish the objectives together and in a correct way.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This week during one of the classes we played a game in which one player used one plastic circle and the other player used 4. The game was that there was a board and the player with 4 circles had to get the 4 pieces to the other side of the board; the same for the player with one piece. I could tell that it seemed that it was easier for the one-piece player to win, however, once the player with 4 pieces learned how to play with the 4 figures, he had a better chance of winning. The way I plan to contribute to my team is by keeping a positive attitude, not being so introverted, but more extroverted so that I can give my suggestions. Likewise, I plan to contribute to the team, helping them to maintain a team atmosphere so we can accomplish the objectives together and in a correct way.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Kevin, a BYUI graduate and a very good programmer, had received a job offer in Tennessee at a company called Slick Tix. During that time, Kevin was maintaining a long-distance relationship with his girlfriend. During his time at work, he was mostly doing sales instead of programming, something that made him uncomfortable at times, as he had been promised something else. At times he felt lonely, as he did not share the same beliefs as his co-workers and did not have a true best friend. Eventually, all these factors led Kevin to have a lot of doubts and to ponder whether he should stay in his current job or quit.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kevin, a BYUI graduate and a very good programmer, had received a job offer in Tennessee at a company called Slick Tix. During that time, Kevin was maintaining a long-distance relationship with his girlfriend. During his time at work, he was mostly doing sales instead of programming, something that made him uncomfortable at times, as he had been promised something else. At times he felt lonely, as he did not share the same beliefs as his co-workers and did not have a true best friend. Eventually, all these factors led Kevin to have a lot of doubts and to ponder whether he should stay in his current job or quit.",
    2) | Out-Null

# --- Add the new trailing content: a blank paragraph, a paragraph summarizing
# --- the number of points, and a "Total Points: 100" paragraph. These are
# --- inserted right after the "I was also able to learn..." paragraph, i.e.
# --- right before the (previously) first trailing empty paragraph. We locate
# --- that empty paragraph, then replace its contents with the 3 new
# --- paragraphs followed by a re-created copy of that same empty paragraph so
# --- none of the original content is lost. ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("I was also able to learn")) {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($target -ne $null) {
    $newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">I answered all the questions making </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>100 points</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> giving my own thoughts of the readings and the teamwork in class.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Total Points: 100</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
    $target.Range.InsertXML($newXml)
}
